# Generate Report for Handoff
#
# The file "daa8f393-dfed-4424-a87a-f9ac622385e9.md" has been (re-)handed
# off for localization, so its status flips from "In Translation" to
# "Ready for handoff" on every sheet that tracks it, and the per-locale
# "Latest Handoff Datetime" timestamps on the locale sheets are refreshed
# to reflect the new handoff.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row for daa8f393-... (row 6, since row 1 is the header)
$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("B6").Value = "Ready for handoff"
$ovw.Range("C6").Value = "Ready for handoff"

# --- zh-cn sheet: status + refreshed handoff datetime
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B6").Value = "Ready for handoff"
$zhcn.Range("D6").Value = "2016-02-22 13:33:45"

# --- de-de sheet: status + refreshed handoff datetime
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B6").Value = "Ready for handoff"
$dede.Range("D6").Value = "2016-02-22 13:33:58"
